$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression: tiny precision change in B2
$ws.Range("B2").Value = 0.3319943219095314

# Row 3 - RandomForestRegressor: values updated
$ws.Range("B3").Value = 0.9852314767677984
$ws.Range("C3").Value = 0.9847658646738863
$ws.Range("D3").Value = 0.8160157250112724

# Row 4 - model name changed from GradientBoostingRegressor to DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9806874230962386
$ws.Range("C4").Value = 0.9798692055097598
$ws.Range("D4").Value = 0.8123113449298193

# Row 5 - model name changed from AdaBoostRegressor to MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.8469489999053322
$ws.Range("C5").Value = 0.8387462818560829
$ws.Range("D5").Value = 0.6465802461876727
